$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 541; existing rows 541:590 shift down to 542:591
$ws.Rows("541:541").Insert()

# Populate the newly inserted row 541 with the new record's data
$ws.Cells.Item(541, 1).Value = 5
$ws.Cells.Item(541, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(541, 3).Value = "Maule"
$ws.Cells.Item(541, 4).Value = 45166
$ws.Cells.Item(541, 5).Value = 7
$ws.Cells.Item(541, 6).Value = 100114013
$ws.Cells.Item(541, 7).Value = "Zanahoria"
$ws.Cells.Item(541, 8).Value = "Sin especificar"
$ws.Cells.Item(541, 9).Value = "Primera"
$ws.Cells.Item(541, 10).Value = 600
$ws.Cells.Item(541, 11).Value = 5500
$ws.Cells.Item(541, 12).Value = 5500
$ws.Cells.Item(541, 13).Value = 5500
$ws.Cells.Item(541, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(541, 15).Value = "Región de Ñuble"
$ws.Cells.Item(541, 16).Value = 275
$ws.Cells.Item(541, 17).Value = 20
$ws.Cells.Item(541, 18).Value = "Hortaliza"
